$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'270.13"
$ws.Range("E2").Value = "'3.22%"
$ws.Range("D3").Value = "'26.70"
$ws.Range("E3").Value = "'-1.46%"
$ws.Range("D4").Value = "'4.713"
$ws.Range("E4").Value = "'0.04%"
$ws.Range("D5").Value = "'0.06109"
$ws.Range("E5").Value = "'-1.58%"
$ws.Range("D6").Value = "'6.739"
$ws.Range("E6").Value = "'0.16%"
$ws.Range("D7").Value = "'0.8566"
$ws.Range("E7").Value = "'0.77%"
$ws.Range("D8").Value = "'0.8981"
$ws.Range("E8").Value = "'-1.43%"
$ws.Range("D9").Value = "'0.1432"
$ws.Range("E9").Value = "'1.68%"
$ws.Range("D10").Value = "'0.04972"
$ws.Range("E10").Value = "'6.12%"
$ws.Range("D11").Value = "'0.07108"
$ws.Range("E11").Value = "'0.27%"
$ws.Range("D12").Value = "'0.03170"
$ws.Range("E12").Value = "'0.40%"
$ws.Range("D13").Value = "'0.09033"
$ws.Range("E13").Value = "'-0.31%"
$ws.Range("D14").Value = "'0.001535"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("D15").Value = "'0.0006081"
$ws.Range("E15").Value = "'-1.09%"
$ws.Range("D16").Value = "'0.006015"
$ws.Range("E16").Value = "'-1.21%"
$ws.Range("E17").Value = "'-0.19%"
$ws.Range("D18").Value = "'3.175"
$ws.Range("E18").Value = "'0.15%"
$ws.Range("E19").Value = "'3.97%"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("D21").Value = "'0.1280"
$ws.Range("E21").Value = "'-1.51%"
$ws.Range("D22").Value = "'3.850"
$ws.Range("E22").Value = "'-5.82%"
$ws.Range("D23").Value = "'0.04246"
$ws.Range("E23").Value = "'0.49%"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("D25").Value = "'0.004152"
$ws.Range("E25").Value = "'0.49%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D27").Value = "'0.0001681"
$ws.Range("E27").Value = "'4.99%"
$ws.Range("D40").Value = "'0.03950"
$ws.Range("E40").Value = "'1.28%"
$ws.Range("E41").Value = "'0.58%"
$ws.Range("D42").Value = "'0.004186"
$ws.Range("E42").Value = "'1.32%"
$ws.Range("D44").Value = "'0.01329"
$ws.Range("E44").Value = "'-1.48%"
$ws.Range("D45").Value = "'0.00005117"
$ws.Range("E45").Value = "'-1.11%"
$ws.Range("E46").Value = "'-0.11%"
$ws.Range("E47").Value = "'-31.85%"
$ws.Range("D48").Value = "'0.9503"
$ws.Range("E48").Value = "'470.20%"
$ws.Range("E49").Value = "'-0.11%"
$ws.Range("E50").Value = "'-0.11%"
